# Weekly fruit/vegetable price update: insert a new daily record for
# "Ciboulette" at Vega Modelo de Temuco, pushing all subsequent records
# down by one row (the series is kept in reverse-chronological insert
# order, newest entry on top of the existing block starting at row 100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: insert a row at 100, shifting rows
# 100-228 down to 101-229 (formatting carries down with the insert).
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A100").Value = 10
$ws.Range("B100").Value = "Vega Modelo de Temuco"
$ws.Range("C100").Value = "La Araucanía"
$ws.Range("D100").Value = 44638
$ws.Range("E100").Value = 9
$ws.Range("F100").Value = 100112039
$ws.Range("G100").Value = "Ciboulette"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 30
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = 5000
$ws.Range("N100").Value = "$/docena de atados"
$ws.Range("O100").Value = "Provincia de Cautín"
$ws.Range("P100").Value = 1667
$ws.Range("Q100").Value = 3
$ws.Range("R100").Value = "Hortaliza"
